$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename table / header columns ---------------------------------------
# Column F (table column 3) was "Dä%&e" -> "Amount"
# Column H (table column 5) was "Date"  -> "Price"
# (updating the header cell text also renames the ListObject's table column)
$ws.Range("F3").Value = "Amount"
$ws.Range("H3").Value = "Price"

# --- Replace the date values in column H with plain numbers --------------
$ws.Range("H4").Style = "Standard"
$ws.Range("H5").Style = "Standard"
$ws.Range("H7").Style = "Standard"
$ws.Range("H9").Style = "Standard"

$ws.Range("H4").Value = 30
$ws.Range("H5").Value = 40
$ws.Range("H7").Value = 60
$ws.Range("H9").Value = 80

# --- Update the selection shown in the sheet view -------------------------
$ws.Range("H10").Select()
